{"js": "const replacements = [\n  [\"57\u00f73=\", \"74\u00f78=\"],\n  [\"24\u00f79=\", \"67\u00f78=\"],\n  [\"47\u00f79=\", \"22\u00f79=\"],\n  [\"26\u00f79=\", \"86\u00f79=\"],\n  [\"33\u00f74=\", \"74\u00f72=\"],\n  [\"15\u00f76=\", \"54\u00f72=\"],\n  [\"88\u00f78=\", \"83\u00f79=\"],\n  [\"12\u00f79=\", \"17\u00f79=\"],\n  [\"54\u00f72=\", \"98\u00f72=\"],\n  [\"54\u00f78=\", \"50\u00f76=\"],\n  [\"77\u00f75=\", \"12\u00f79=\"],\n  [\"62\u00f72=\", \"10\u00f79=\"],\n  [\"78\u00f75=\", \"36\u00f79=\"],\n  [\"64\u00f75=\", \"29\u00f72=\"],\n  [\"60\u00f77=\", \"28\u00f73=\"],\n  [\"91\u00f75=\", \"63\u00f72=\"],\n  [\"85\u00f78=\", \"58\u00f76=\"],\n  [\"12\u00f77=\", \"29\u00f78=\"],\n  [\"44\u00f72=\", \"92\u00f72=\"],\n  [\"76\u00f73=\", \"93\u00f73=\"],\n  [\"79\u00f72=\", \"45\u00f79=\"],\n  [\"30\u00f75=\", \"34\u00f77=\"],\n  [\"38\u00f75=\", \"15\u00f78=\"],\n  [\"83\u00f79=\", \"38\u00f74=\"],\n  [\"43\u00f74=\", \"13\u00f74=\"],\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Snapshot each paragraph's current text first (one sync), then decide\n// and apply replacements (second sync) \u2014 this avoids a later rule's\n// replacement text ever being mistaken for a still-unprocessed \"before\"\n// value, since every match decision is made against the original text.\nconst ranges = paragraphs.items.map((p) => p.getRange());\nfor (const range of ranges) {\n  range.load(\"text\");\n}\nawait context.sync();\n\nfor (const range of ranges) {\n  const original = range.text;\n  for (const [before, after] of replacements) {\n    if (original === before) {\n      range.insertText(after, Word.InsertLocation.replace);\n      break;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"57\u00f73=\", \"74\u00f78=\"),\n    @(\"24\u00f79=\", \"67\u00f78=\"),\n    @(\"47\u00f79=\", \"22\u00f79=\"),\n    @(\"26\u00f79=\", \"86\u00f79=\"),\n    @(\"33\u00f74=\", \"74\u00f72=\"),\n    @(\"15\u00f76=\", \"54\u00f72=\"),\n    @(\"88\u00f78=\", \"83\u00f79=\"),\n    @(\"12\u00f79=\", \"17\u00f79=\"),\n    @(\"54\u00f72=\", \"98\u00f72=\"),\n    @(\"54\u00f78=\", \"50\u00f76=\"),\n    @(\"77\u00f75=\", \"12\u00f79=\"),\n    @(\"62\u00f72=\", \"10\u00f79=\"),\n    @(\"78\u00f75=\", \"36\u00f79=\"),\n    @(\"64\u00f75=\", \"29\u00f72=\"),\n    @(\"60\u00f77=\", \"28\u00f73=\"),\n    @(\"91\u00f75=\", \"63\u00f72=\"),\n    @(\"85\u00f78=\", \"58\u00f76=\"),\n    @(\"12\u00f77=\", \"29\u00f78=\"),\n    @(\"44\u00f72=\", \"92\u00f72=\"),\n    @(\"76\u00f73=\", \"93\u00f73=\"),\n    @(\"79\u00f72=\", \"45\u00f79=\"),\n    @(\"30\u00f75=\", \"34\u00f77=\"),\n    @(\"38\u00f75=\", \"15\u00f78=\"),\n    @(\"83\u00f79=\", \"38\u00f74=\"),\n    @(\"43\u00f74=\", \"13\u00f74=\")\n)\n\n# Two-phase replacement to avoid a later rule's target text accidentally\n# matching (and re-mutating) a cell already rewritten by an earlier rule:\n#   phase 1: before-text -> unique placeholder token (by index)\n#   phase 2: placeholder token -> final after-text\n\nfor ($i = 0; $i -lt $replacements.Count; $i++) {\n    $before = $replacements[$i][0]\n    $token = \"@@TOKEN\" + $i + \"@@\"\n\n    foreach ($tbl in $d.Tables) {\n        foreach ($cell in $tbl.Range.Cells) {\n            $cellRange = $cell.Range\n            $find = $cellRange.Find\n            $find.ClearFormatting()\n            $find.Replacement.ClearFormatting()\n            $find.Execute(\n                $before,      # FindText\n                $false,       # MatchCase\n                $true,        # MatchWholeWord\n                $false,       # MatchWildcards\n                $false,       # MatchSoundsLike\n                $false,       # MatchAllWordForms\n                $true,        # Forward\n                1,            # Wrap (wdFindContinue)\n                $false,       # Format\n                $token,       # ReplaceWith\n                2             # Replace (wdReplaceOne)\n            ) | Out-Null\n        }\n    }\n}\n\nfor ($i = 0; $i -lt $replacements.Count; $i++) {\n    $after = $replacements[$i][1]\n    $token = \"@@TOKEN\" + $i + \"@@\"\n\n    foreach ($tbl in $d.Tables) {\n        foreach ($cell in $tbl.Range.Cells) {\n            $cellRange = $cell.Range\n            $find = $cellRange.Find\n            $find.ClearFormatting()\n            $find.Replacement.ClearFormatting()\n            $find.Execute(\n                $token,       # FindText\n                $false,       # MatchCase\n                $true,        # MatchWholeWord\n                $false,       # MatchWildcards\n                $false,       # MatchSoundsLike\n                $false,       # MatchAllWordForms\n                $true,        # Forward\n                1,            # Wrap (wdFindContinue)\n                $false,       # Format\n                $after,       # ReplaceWith\n                2             # Replace (wdReplaceOne)\n            ) | Out-Null\n        }\n    }\n}\n"}
